$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1180.0588
$ws.Range("I28").Value = 380
$ws.Range("K28").Value = 380
$ws.Range("M28").Value = 105
$ws.Range("H62").Value = 5040
$ws.Range("I62").Value = 1343.4
$ws.Range("K62").Value = 1343.4
$ws.Range("M62").Value = -719.4000000000001
$ws.Range("H64").Value = 7471.517
$ws.Range("I64").Value = 2997.5
$ws.Range("K64").Value = 2997.5
$ws.Range("M64").Value = -2749.5
$ws.Range("H65").Value = 5040
$ws.Range("I65").Value = 1343.4
$ws.Range("K65").Value = 6717
$ws.Range("M65").Value = -3597
$ws.Range("H67").Value = 7471.517
$ws.Range("I67").Value = 2997.5
$ws.Range("K67").Value = 2997.5
$ws.Range("M67").Value = -2139.5
$ws.Range("H76").Value = 6703
$ws.Range("I76").Value = 5691.5
$ws.Range("J76").Value = 7570
$ws.Range("K76").Value = 5691.5
$ws.Range("L76").Value = 7570
$ws.Range("M76").Value = -5376.5
$ws.Range("N76").Value = -8200
$ws.Range("H79").Value = 6703
$ws.Range("I79").Value = 5691.5
$ws.Range("J79").Value = 7570
$ws.Range("K79").Value = 5691.5
$ws.Range("L79").Value = 7570
$ws.Range("M79").Value = -4599.5
$ws.Range("N79").Value = -9754
$ws.Range("H86").Value = 9099.333000000001
$ws.Range("I86").Value = 10915
$ws.Range("J86").Value = 7888.8887
$ws.Range("K86").Value = 10915
$ws.Range("L86").Value = 7888.8887
$ws.Range("M86").Value = -9792
$ws.Range("N86").Value = -10134.8887
$ws.Range("H89").Value = 9099.333000000001
$ws.Range("I89").Value = 10915
$ws.Range("J89").Value = 7888.8887
$ws.Range("K89").Value = 54575
$ws.Range("L89").Value = 39444.4435
$ws.Range("M89").Value = -48959
$ws.Range("N89").Value = -50676.4435
$ws.Range("H92").Value = 2211.348
$ws.Range("I92").Value = 1587.375
$ws.Range("J92").Value = 3637.5715
$ws.Range("K92").Value = 1587.375
$ws.Range("L92").Value = 3637.5715
$ws.Range("M92").Value = -339.375
$ws.Range("N92").Value = -6133.5715
$ws.Range("H100").Value = 634.3823
$ws.Range("I100").Value = 572.7143
$ws.Range("K100").Value = 572.7143
$ws.Range("M100").Value = -31.71429999999998
$ws.Range("H107").Value = 37148492
$ws.Range("I107").Value = 55555908
$ws.Range("K107").Value = 55555908
$ws.Range("M107").Value = -55553988
$ws.Range("H116").Value = 5332.207
$ws.Range("I116").Value = 3604.7646
$ws.Range("J116").Value = 7779.4165
$ws.Range("K116").Value = 3604.7646
$ws.Range("L116").Value = 7779.4165
$ws.Range("M116").Value = -162.7646
$ws.Range("N116").Value = -14663.4165
$ws.Range("H121").Value = 2276.5386
$ws.Range("J121").Value = 2276.5386
$ws.Range("L121").Value = 6829.6158
$ws.Range("N121").Value = -10323.6158
$ws.Range("H132").Value = 4938
$ws.Range("I132").Value = 5590.9443
$ws.Range("K132").Value = 16772.8329
$ws.Range("M132").Value = -14242.8329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1286472.9
$ws.Range("I2").Value = 1414820.9
$ws.Range("J2").Value = 2993.5
$ws.Range("K2").Value = 1414820.9
$ws.Range("L2").Value = 2993.5
$ws.Range("M2").Value = -1414707.9
$ws.Range("N2").Value = -3219.5
$ws.Range("H32").Value = 10022.373
$ws.Range("I32").Value = 6114.1953
$ws.Range("J32").Value = 18924.334
$ws.Range("K32").Value = 6114.1953
$ws.Range("L32").Value = 18924.334
$ws.Range("M32").Value = -5827.1953
$ws.Range("N32").Value = -19498.334
$ws.Range("H61").Value = 6587
$ws.Range("I61").Value = 6746.2
$ws.Range("K61").Value = 6746.2
$ws.Range("M61").Value = -6534.2
$ws.Range("H97").Value = 1526123.8
$ws.Range("I97").Value = 1598772.5
$ws.Range("K97").Value = 1598772.5
$ws.Range("M97").Value = -1598276.5
$ws.Range("H102").Value = 2527055.5
$ws.Range("I102").Value = 2605869.8
$ws.Range("K102").Value = 2605869.8
$ws.Range("M102").Value = -2604247.8
$ws.Range("H116").Value = 1286472.9
$ws.Range("I116").Value = 1414820.9
$ws.Range("J116").Value = 2993.5
$ws.Range("K116").Value = 1414820.9
$ws.Range("L116").Value = 2993.5
$ws.Range("M116").Value = -1412526.9
$ws.Range("N116").Value = -7581.5
$ws.Range("H118").Value = 56166.332
$ws.Range("J118").Value = 56166.332
$ws.Range("L118").Value = 56166.332
$ws.Range("N118").Value = -59480.332
$ws.Range("H132").Value = 1914.6136
$ws.Range("I132").Value = 1127.5714
$ws.Range("K132").Value = 3382.7142
$ws.Range("M132").Value = -852.7142000000003
$ws.Range("H136").Value = 6587
$ws.Range("I136").Value = 6746.2
$ws.Range("K136").Value = 20238.6
$ws.Range("M136").Value = -17688.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1286472.9
$ws.Range("I3").Value = 1414820.9
$ws.Range("J3").Value = 2993.5
$ws.Range("K3").Value = 1414820.9
$ws.Range("L3").Value = 2993.5
$ws.Range("M3").Value = -1414706.9
$ws.Range("N3").Value = -3221.5
$ws.Range("H56").Value = 30000
$ws.Range("J56").Value = 30000
$ws.Range("L56").Value = 30000
$ws.Range("N56").Value = -31478
$ws.Range("H94").Value = 2946888
$ws.Range("I94").Value = 3226877.2
$ws.Range("J94").Value = 53666.668
$ws.Range("K94").Value = 3226877.2
$ws.Range("L94").Value = 53666.668
$ws.Range("M94").Value = -3226426.2
$ws.Range("N94").Value = -54568.668
$ws.Range("H134").Value = 6402.476
$ws.Range("I134").Value = 1555.3684
$ws.Range("K134").Value = 4666.1052
$ws.Range("M134").Value = -2131.1052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 818.75
$ws.Range("I105").Value = 545
$ws.Range("K105").Value = 545
$ws.Range("M105").Value = 1202
$ws.Range("H132").Value = 88422.05
$ws.Range("I132").Value = 57827.11
$ws.Range("K132").Value = 173481.33
$ws.Range("M132").Value = -170951.33
$ws.Range("H134").Value = 2025.4117
$ws.Range("I134").Value = 1352.6428
$ws.Range("J134").Value = 5165
$ws.Range("K134").Value = 4057.9284
$ws.Range("L134").Value = 15495
$ws.Range("M134").Value = -1522.9284
$ws.Range("N134").Value = -20565

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 5100
$ws.Range("I60").Value = 200
$ws.Range("K60").Value = 600
$ws.Range("M60").Value = -349
$ws.Range("H61").Value = 201.8
$ws.Range("I61").Value = 52
$ws.Range("J61").Value = 301.66666
$ws.Range("K61").Value = 156
$ws.Range("L61").Value = 904.9999799999999
$ws.Range("M61").Value = 59
$ws.Range("N61").Value = -1334.99998
$ws.Range("H81").Value = 6240.48
$ws.Range("I81").Value = 763.5
$ws.Range("J81").Value = 7283.7144
$ws.Range("K81").Value = 2290.5
$ws.Range("L81").Value = 21851.1432
$ws.Range("M81").Value = -1167.5
$ws.Range("N81").Value = -24097.1432
$ws.Range("H84").Value = 6240.48
$ws.Range("I84").Value = 763.5
$ws.Range("J84").Value = 7283.7144
$ws.Range("K84").Value = 6871.5
$ws.Range("L84").Value = 65553.4296
$ws.Range("M84").Value = -1255.5
$ws.Range("N84").Value = -76785.4296
$ws.Range("H93").Value = 5050
$ws.Range("I93").Value = 10000
$ws.Range("J93").Value = 100
$ws.Range("K93").Value = 30000
$ws.Range("L93").Value = 300
$ws.Range("M93").Value = -28128
$ws.Range("N93").Value = -4044
$ws.Range("H137").Value = 2327.7693
$ws.Range("I137").Value = 1893.2
$ws.Range("J137").Value = 3776.3333
$ws.Range("K137").Value = 5679.6
$ws.Range("L137").Value = 11328.9999
$ws.Range("M137").Value = -579.6000000000004
$ws.Range("N137").Value = -21528.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3392415.5
$ws.Range("I126").Value = 2275129.8
$ws.Range("K126").Value = 6825389.399999999
$ws.Range("M126").Value = -6822919.399999999
$ws.Range("H132").Value = 2208.228
$ws.Range("I132").Value = 1966.92
$ws.Range("K132").Value = 5900.76
$ws.Range("M132").Value = -3370.76
$ws.Range("H141").Value = 51836.57
$ws.Range("J141").Value = 71964
$ws.Range("L141").Value = 71964
$ws.Range("N141").Value = -82324

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 37889.12
$ws.Range("I22").Value = 74819.914
$ws.Range("J22").Value = 3799.1538
$ws.Range("K22").Value = 74819.914
$ws.Range("L22").Value = 3799.1538
$ws.Range("M22").Value = -74524.914
$ws.Range("N22").Value = -4389.1538
$ws.Range("H27").Value = 37889.12
$ws.Range("I27").Value = 74819.914
$ws.Range("J27").Value = 3799.1538
$ws.Range("K27").Value = 74819.914
$ws.Range("L27").Value = 3799.1538
$ws.Range("M27").Value = -74712.914
$ws.Range("N27").Value = -4013.1538
$ws.Range("H88").Value = 14647.5
$ws.Range("I88").Value = 11539.2
$ws.Range("J88").Value = 30189
$ws.Range("K88").Value = 11539.2
$ws.Range("L88").Value = 30189
$ws.Range("M88").Value = -11111.2
$ws.Range("N88").Value = -31045
$ws.Range("H91").Value = 14647.5
$ws.Range("I91").Value = 11539.2
$ws.Range("J91").Value = 30189
$ws.Range("K91").Value = 11539.2
$ws.Range("L91").Value = 30189
$ws.Range("M91").Value = -10057.2
$ws.Range("N91").Value = -33153
$ws.Range("H132").Value = 6638.392
$ws.Range("I132").Value = 6773.548
$ws.Range("J132").Value = 6007.6665
$ws.Range("K132").Value = 20320.644
$ws.Range("L132").Value = 18022.9995
$ws.Range("M132").Value = -17790.644
$ws.Range("N132").Value = -23082.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 43479836
$ws.Range("I107").Value = 100000680
$ws.Range("J107").Value = 2262.4614
$ws.Range("K107").Value = 300002040
$ws.Range("L107").Value = 6787.3842
$ws.Range("M107").Value = -300000120
$ws.Range("N107").Value = -10627.3842
$ws.Range("H132").Value = 20221932
$ws.Range("I132").Value = 27782742
$ws.Range("J132").Value = 779850.9399999999
$ws.Range("K132").Value = 83348226
$ws.Range("L132").Value = 2339552.82
$ws.Range("M132").Value = -83345696
$ws.Range("N132").Value = -2344612.82
$ws.Range("H136").Value = 1036.8064
$ws.Range("I136").Value = 840.21155
$ws.Range("K136").Value = 2520.63465
$ws.Range("M136").Value = 29.36535000000003
